$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.02530399975908262
$ws.Range("E2").Value = 0.02530399975908262

# Row 3
$ws.Range("D3").Value = 0.8738244904266079
$ws.Range("E3").Value = 0.8738244904266079

# Row 4
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = 0.5263704196157762
$ws.Range("E4").Value = 0.5263704196157762

# Row 5
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = 0.08263331213056946
$ws.Range("E5").Value = 0.08263331213056946

# Row 6
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = 0.04332418076429817
$ws.Range("E6").Value = 0.04332418076429817

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.06099346724616209
$ws.Range("E7").Value = 0.9390065327538379

# Row 8
$ws.Range("D8").Value = 0.01619906853771051
$ws.Range("E8").Value = 0.9838009314622895

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.2077240461859528
$ws.Range("E9").Value = 0.7922759538140471

# Row 10
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = 0.5341682101457098
$ws.Range("E10").Value = 0.4658317898542902
$ws.Range("F10").Value = 1.34355103969574

# Row 11
$ws.Range("D11").Value = 0.06147103592105861
$ws.Range("E11").Value = 0.06147103592105861

# Row 12
$ws.Range("D12").Value = 0.9557763857290811
$ws.Range("E12").Value = 0.9557763857290811

# Row 13
$ws.Range("C13").Value = $false
$ws.Range("D13").Value = 0.5150759946714326
$ws.Range("E13").Value = 0.5150759946714326

# Row 14
$ws.Range("C14").Value = $true
$ws.Range("D14").Value = 0.1244465543032987
$ws.Range("E14").Value = 0.1244465543032987

# Row 15
$ws.Range("C15").Value = $true
$ws.Range("D15").Value = 0.02304341319258804
$ws.Range("E15").Value = 0.02304341319258804

# Row 16
$ws.Range("C16").Value = $false
$ws.Range("D16").Value = 0.06383976400078317
$ws.Range("E16").Value = 0.9361602359992168

# Row 17
$ws.Range("D17").Value = 0.008592741989594525
$ws.Range("E17").Value = 0.9914072580104055

# Row 18
$ws.Range("D18").Value = 0.5101654295631826
$ws.Range("E18").Value = 0.4898345704368174

# Row 19
$ws.Range("D19").Value = 0.7844383574811579
$ws.Range("E19").Value = 0.2155616425188421
$ws.Range("F19").Value = 1.38732647895813
